$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.297.95"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.834.74"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.82"
$ws.Range("E5").Value = "  -5.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.50"
$ws.Range("E6").Value = "  -8.91%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  -5.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.832.19"
$ws.Range("E9").Value = "  -4.63%  "
$ws.Range("E10").Value = "  -7.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.93"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.326.37"
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.376.08"
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.81"
$ws.Range("E16").Value = "  -7.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.840.48"
$ws.Range("E17").Value = "  -4.57%  "
$ws.Range("E18").Value = "  -6.54%  "
$ws.Range("E19").Value = "  -7.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.08"
$ws.Range("E20").Value = "  -6.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.56"
$ws.Range("E21").Value = "  -5.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("E22").Value = "  -6.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.15"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("E26").Value = "  -7.78%  "
$ws.Range("E27").Value = "  -7.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  -7.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0813"
$ws.Range("E30").Value = "  -9.92%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.02"
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.70"
$ws.Range("E34").Value = "  -5.47%  "
$ws.Range("E35").Value = "  -7.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.37"
$ws.Range("E36").Value = "  -7.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.926"
$ws.Range("E37").Value = "  -11.22%  "
$ws.Range("E38").Value = "  -8.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.56"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.220.84"
$ws.Range("E40").Value = "  -7.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.629"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").Value = "  -7.95%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.37"
$ws.Range("E43").Value = "  -9.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0561"
$ws.Range("E44").Value = "  -4.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.49"
$ws.Range("E46").Value = "  -10.09%  "
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0226"
$ws.Range("E48").Value = "  -6.19%  "
$ws.Range("E49").Value = "  -5.70%  "
$ws.Range("E50").Value = "  -12.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.72"
$ws.Range("E51").Value = "  -7.55%  "
